$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.639.06"
$ws.Range("D3").Value = "'1.592.45"
$ws.Range("E3").Value = "  -1.67%  "
$ws.Range("D5").Value = "'212.10"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -1.46%  "
$ws.Range("D9").Value = "'0.245"
$ws.Range("E9").Value = "  -2.81%  "
$ws.Range("D10").Value = "'19.60"
$ws.Range("E10").Value = "  -2.20%  "
$ws.Range("D11").Value = "'0.0835"
$ws.Range("E11").Value = "  -1.56%  "
$ws.Range("D12").Value = "'1.815.99"
$ws.Range("E12").Value = "  -1.65%  "
$ws.Range("D13").Value = "'1.590.52"
$ws.Range("E13").Value = "  -1.86%  "
$ws.Range("E14").Value = "  -2.80%  "
$ws.Range("E15").Value = "  -2.90%  "
$ws.Range("D16").Value = "'65.16"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").Value = "'26.608.73"
$ws.Range("E17").Value = "  -1.46%  "
$ws.Range("D18").Value = "'0.0₃0729"
$ws.Range("E18").Value = "  -2.46%  "
$ws.Range("E19").Value = "  -2.68%  "
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "'6.69"
$ws.Range("E21").Value = "  -2.39%  "
$ws.Range("D22").Value = "'4.26"
$ws.Range("E22").Value = "  -2.26%  "
$ws.Range("E23").Value = "  -3.00%  "
$ws.Range("E24").Value = "  -1.76%  "
$ws.Range("D25").Value = "'145.92"
$ws.Range("E25").Value = "  -1.59%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").Value = "'7.15"
$ws.Range("E27").Value = "  -3.67%  "
$ws.Range("E28").Value = "  -1.76%  "
$ws.Range("E29").Value = "  -1.63%  "
$ws.Range("D30").Value = "'0.0504"
$ws.Range("E30").Value = "  -2.03%  "
$ws.Range("E31").Value = "  -1.51%  "
$ws.Range("E32").Value = "  -3.89%  "
$ws.Range("D33").Value = "'0.667"
$ws.Range("E33").Value = "  -10.89%  "
$ws.Range("E34").Value = "  -3.42%  "
$ws.Range("D35").Value = "'1.298.90"
$ws.Range("E35").Value = "  -3.59%  "
$ws.Range("E36").Value = "  -0.42%  "
$ws.Range("E37").Value = "  -5.02%  "
$ws.Range("E38").Value = "  -3.85%  "
$ws.Range("E39").Value = "  -2.31%  "
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("E41").Value = "  -1.06%  "
$ws.Range("E42").Value = "  +0.69%  "
$ws.Range("D43").Value = "'2.19"
$ws.Range("E43").Value = "  -1.94%  "
$ws.Range("D44").Value = "'63.15"
$ws.Range("E44").Value = "  -3.27%  "
$ws.Range("D45").Value = "'1.728.75"
$ws.Range("E45").Value = "  -1.55%  "
$ws.Range("D46").Value = "'89.03"
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("D48").Value = "'0.800"
$ws.Range("E48").Value = "  -7.55%  "
$ws.Range("D49").Value = "'0.0983"
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("D50").Value = "'0.0502"
$ws.Range("E50").Value = "  -2.72%  "
$ws.Range("D51").Value = "'7.55"
$ws.Range("E51").Value = "  -1.88%  "
